$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Jessica row (row 2) ---
$ws.Range("B2").Value = "Thinks she is brave going around during the conflict. A bit too free."
$ws.Range("C2").Value = "Sees as being friendly and dad-like."
$ws.Range("D2").Value = "For once does not trust someone at the door. Questions their authority but is quiet/shy/scared."
$ws.Range("E2").Value = "Trusting. Clings to him a bit, seeing her own dad in him. Makes her think about how her parents left her."
$ws.Range("F2").Value = "Doesn't want to judge a book by it's cover, but is terrified once let in."
$ws.Range("J2").ClearContents()

# --- Bob row (row 3) ---
$ws.Range("B3").Value = "Thinks she is a bit bonkers and wondering if there is a hidden agenda."
$ws.Range("C3").Value = "Neutral. Wary, but leaning towards trustworthy."
$ws.Range("D3").Value = "Readys for a fight. Thinks there are no possible good intentions."
$ws.Range("E3").Value = "Appreciates how he looks out for his son, but distrusts how he asks for the gun."
$ws.Range("F3").Value = "Doesn't want to let in a massive, scratching shape. Will be hostile towards it if let in."

# --- Violet row (row 4) ---
$ws.Range("C4").Value = "Neutral. He's just another guy, could be helpful."
$ws.Range("D4").Value = "Similar to Bob. Adrenaline starts pumping and she gets ready to defend."
$ws.Range("E4").Value = "Relatively neutral. Thinks his mission is fruitless but doesn't tell him that."

# --- Hal row (row 5) ---
$ws.Range("D5").Value = "Hateful/angry. These kinds of people woulded Sal."

# --- Sal row (row 6) ---
$ws.Range("D6").Value = "Angry/nervous of harm coming to Hal. Tries to remain stoic and calm."

# --- Dad row (row 8) ---
$ws.Range("C8").Value = "Sees as a fellow human being and a good guy. Appreciates beard."
$ws.Range("K8").Value = "Empathizes with her struggles and loss of parents. She motivates him even further to find his son."
$ws.Range("J8").Value = "Doesn't talk to Bob much. Realizes that Bob is a bit distrustful and chooses to ignore him."
$ws.Range("I8").Value = "Neutral. Likes her resilience and hands-on nature."

# --- Raiders row (row 9) ---
$ws.Range("G9").Value = "Kill. See fighting-age male as hostile."
$ws.Range("H9").Value = "Kidnap."
$ws.Range("I9").Value = "Kidnap."
$ws.Range("J9").Value = "Kill."
$ws.Range("K9").Value = "Kill. See disability as a hassle to deal with."

# --- Neighbour row (row 10) ---
$ws.Range("D10").Value = "Untrusting. Finds them suspicious."

# Cells that become "N/A" (re-using the existing grey N/A style already used elsewhere,
# e.g. F4) via a format-only paste, then writing the N/A text.
$ws.Range("F4").Copy()
$naCells = @("E10", "I10", "J10", "K10", "J11", "K11")
foreach ($cellRef in $naCells) {
    $ws.Range($cellRef).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
foreach ($cellRef in $naCells) {
    $ws.Range($cellRef).Value = "N/A"
}

# --- Selection moved by the author while reviewing the finished sheet ---
$ws.Range("I6").Select()
